# Weekly update to "Hortaliza, Vega Modelo de Temuco - Arveja Verde":
# a new price-report row is inserted at row 46 (pushing the existing
# rows 46-136 down to 47-137), and populated with the latest record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46; Excel shifts rows 46:136 -> 47:137
# and grows the used range to A1:R137 (matching row height/format of the
# row being pushed down, same as a normal Excel "Insert Row" above D46).
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new observation.
$ws.Cells.Item(46, 1).Value  = 10
$ws.Cells.Item(46, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(46, 3).Value  = "La Araucanía"
$ws.Cells.Item(46, 4).Value  = 45272
$ws.Cells.Item(46, 5).Value  = 9
$ws.Cells.Item(46, 6).Value  = 100112022
$ws.Cells.Item(46, 7).Value  = "Arveja Verde"
$ws.Cells.Item(46, 8).Value  = "Sin especificar"
$ws.Cells.Item(46, 9).Value  = "Primera"
$ws.Cells.Item(46, 10).Value = 55
$ws.Cells.Item(46, 11).Value = 25000
$ws.Cells.Item(46, 12).Value = 25000
$ws.Cells.Item(46, 13).Value = 25000
$ws.Cells.Item(46, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(46, 15).Value = "Región del Maule"
$ws.Cells.Item(46, 16).Value = 1000
$ws.Cells.Item(46, 17).Value = 25
$ws.Cells.Item(46, 18).Value = "Hortaliza"
